$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "2.0-Harina  (kg),1.0-Huevos (unidad),1.0-Vainilla (ml),3.0-Leche (litros),"
$ws.Range("C4").Value = "5.0-Harina  (kg),2.0-Huevos (unidad),1.0-Vainilla (ml),"
$ws.Range("C6").Value = "4.0-Harina  (kg),2.0-Limon (unidad),5.0-Huevos (unidad),5.0-Crema (litros),"
$ws.Range("C7").Value = "0.3-Harina  (kg),2.0-Huevos (unidad),0.1-Vainilla (ml),0.2-Leche (litros),"
